$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text (not numeric) storage for price/volume columns, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

# Updated coin price/volume snapshot (GitHub Actions symbol-list refresh).
$ws.Range("D2").Value = "314.69"
$ws.Range("E2").Value = "2.47%"
$ws.Range("D3").Value = "39.48"
$ws.Range("E3").Value = "2.38%"
$ws.Range("D4").Value = "5.132"
$ws.Range("E4").Value = "0.67%"
$ws.Range("D5").Value = "0.08173"
$ws.Range("D6").Value = "1.966"
$ws.Range("E6").Value = "0.30%"
$ws.Range("D7").Value = "8.202"
$ws.Range("E7").Value = "3.18%"
$ws.Range("D8").Value = "0.9275"
$ws.Range("E8").Value = "-0.21%"
$ws.Range("E9").Value = "-0.82%"
$ws.Range("D10").Value = "0.1977"
$ws.Range("E10").Value = "0.91%"
$ws.Range("D11").Value = "0.09073"
$ws.Range("E11").Value = "0.62%"
$ws.Range("D12").Value = "0.03496"
$ws.Range("E12").Value = "-0.21%"
$ws.Range("D13").Value = "0.09814"
$ws.Range("E13").Value = "-0.10%"
$ws.Range("D14").Value = "0.001397"
$ws.Range("E14").Value = "-0.65%"
$ws.Range("D15").Value = "0.006116"
$ws.Range("E15").Value = "1.20%"
$ws.Range("D16").Value = "3.661"
$ws.Range("E16").Value = "-1.87%"
$ws.Range("D17").Value = "4.237"
$ws.Range("E17").Value = "0.93%"
$ws.Range("D18").Value = "3.171"
$ws.Range("E18").Value = "-7.25%"
$ws.Range("D19").Value = "0.3457"
$ws.Range("E19").Value = "-0.14%"
$ws.Range("E20").Value = "0.26%"
$ws.Range("D21").Value = "4.746"
$ws.Range("E21").Value = "-1.00%"
$ws.Range("D22").Value = "0.2426"
$ws.Range("E22").Value = "-1.12%"
$ws.Range("D23").Value = "0.04380"
$ws.Range("E23").Value = "-0.56%"
$ws.Range("D24").Value = "0.001222"
$ws.Range("E24").Value = "-0.31%"
$ws.Range("D25").Value = "0.004786"
$ws.Range("E25").Value = "-1.00%"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").Value = "-0.22%"
$ws.Range("D27").Value = "0.0003995"
$ws.Range("E27").Value = "-10.17%"
$ws.Range("D39").Value = "0.02167"
$ws.Range("E39").Value = "4.13%"
$ws.Range("D40").Value = "0.05203"
$ws.Range("E40").Value = "1.87%"
$ws.Range("D41").Value = "0.007518"
$ws.Range("E41").Value = "0.49%"
$ws.Range("D42").Value = "0.009789"
$ws.Range("E42").Value = "-3.48%"
$ws.Range("D43").Value = "0.1371"
$ws.Range("E43").Value = "1.11%"
$ws.Range("D44").Value = "0.002114"
$ws.Range("E44").Value = "-0.91%"
$ws.Range("D45").Value = "0.009129"
$ws.Range("E45").Value = "-1.58%"
$ws.Range("D46").Value = "0.00006388"
$ws.Range("E46").Value = "2.32%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.30%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "0.001199"
$ws.Range("E48").Value = "-25.13%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "0.002766"
$ws.Range("E49").Value = "-8.70%"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").Value = "-0.30%"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").Value = "-0.30%"
